$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price (D) cells whose new values look numeric,
# so Excel does not silently convert them to numbers (they must remain text).
$ws.Range("D4,D5,D6,D9,D11,D12,D14,D16,D19,D20,D21,D22,D23,D24,D26,D28,D29,D30,D31,D32,D33,D36,D37,D38,D39,D40,D41,D42,D43,D45,D48,D49,D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "79.553.71"
$ws.Range("E2").Value = "  +4.23%  "

# Row 3
$ws.Range("D3").Value = "3.167.59"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "206.10"
$ws.Range("E5").Value = "  +4.03%  "

# Row 6
$ws.Range("D6").Value = "627.02"
$ws.Range("E6").Value = "  +1.87%  "

# Row 7
$ws.Range("E7").Value = "  +28.15%  "

# Row 9
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  +7.34%  "

# Row 10
$ws.Range("D10").Value = "3.165.25"
$ws.Range("E10").Value = "  +2.55%  "

# Row 11
$ws.Range("D11").Value = "0.594"
$ws.Range("E11").Value = "  +34.81%  "

# Row 12
$ws.Range("D12").Value = "0.0000253"
$ws.Range("E12").Value = "  +30.62%  "

# Row 13
$ws.Range("E13").Value = "  +2.08%  "

# Row 14
$ws.Range("D14").Value = "5.29"
$ws.Range("E14").Value = "  +1.21%  "

# Row 15
$ws.Range("D15").Value = "3.747.04"
$ws.Range("E15").Value = "  +2.59%  "

# Row 16
$ws.Range("D16").Value = "31.55"
$ws.Range("E16").Value = "  +8.08%  "

# Row 17
$ws.Range("D17").Value = "79.539.59"
$ws.Range("E17").Value = "  +4.36%  "

# Row 18
$ws.Range("D18").Value = "3.156.72"
$ws.Range("E18").Value = "  +2.64%  "

# Row 19
$ws.Range("D19").Value = "14.34"
$ws.Range("E19").Value = "  +5.55%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "438.09"
$ws.Range("E20").Value = "  +14.89%  "

# Row 21
$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").Value = "2.94"
$ws.Range("E21").Value = "  +14.98%  "

# Row 22
$ws.Range("D22").Value = "9.15"
$ws.Range("E22").Value = "  +0.37%  "

# Row 23
$ws.Range("D23").Value = "5.26"
$ws.Range("E23").Value = "  +18.72%  "

# Row 24
$ws.Range("D24").Value = "6.67"
$ws.Range("E24").Value = "  +2.97%  "

# Row 25
$ws.Range("D25").Value = "3.330.77"
$ws.Range("E25").Value = "  +2.85%  "

# Row 26
$ws.Range("D26").Value = "76.08"
$ws.Range("E26").Value = "  +5.20%  "

# Row 27
$ws.Range("E27").Value = "  +5.70%  "

# Row 28
$ws.Range("D28").Value = "10.90"
$ws.Range("E28").Value = "  +9.11%  "

# Row 29
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.37%  "

# Row 30
$ws.Range("D30").Value = "0.0000122"
$ws.Range("E30").Value = "  +12.23%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "9.04"
$ws.Range("E31").Value = "  +8.68%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").Value = "550.85"
$ws.Range("E33").Value = "  +10.10%  "

# Row 34
$ws.Range("E34").Value = "  +3.97%  "

# Row 35
$ws.Range("E35").Value = "  +4.14%  "

# Row 36
$ws.Range("D36").Value = "0.150"
$ws.Range("E36").Value = "  +21.71%  "

# Row 37
$ws.Range("D37").Value = "23.12"
$ws.Range("E37").Value = "  +11.13%  "

# Row 38
$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  +18.41%  "

# Row 39
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("D40").Value = "0.408"
$ws.Range("E40").Value = "  +7.61%  "

# Row 41
$ws.Range("D41").Value = "20.76"
$ws.Range("E41").Value = "  +3.48%  "

# Row 42
$ws.Range("D42").Value = "164.08"
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  +10.12%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "187.88"
$ws.Range("E45").Value = "  -3.78%  "

# Row 46
$ws.Range("E46").Value = "  +9.42%  "

# Row 47
$ws.Range("E47").Value = "  +10.25%  "

# Row 48
$ws.Range("D48").Value = "0.782"
$ws.Range("E48").Value = "  -2.55%  "

# Row 49
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  +4.52%  "

# Row 50
$ws.Range("E50").Value = "  +4.88%  "

# Row 51
$ws.Range("D51").Value = "4.28"
$ws.Range("E51").Value = "  +9.57%  "
